$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15: fill in the 8MHz crystal part details (quantity, description, part number, digikey part number)
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = "CRYSTAL 8.0000MHZ 12PF SMD"
$ws.Range("D15").Value = "CX5032GB08000H0HPQZ1"
$ws.Range("E15").Value = "1253-1373-2-ND"

# New Row 18: through-hole debug header
$ws.Range("B18").Value = 0
$ws.Range("E18").Value = "609-3714-ND"
$ws.Range("G18").Value = "Through hole (across board edge) 2x7 debug header (must make boarod 0.8mm thick)"
$ws.Range("D18").Value = "20021111-00014T4LF"
$ws.Range("C18").Value = "CONN HEADER VERT 14POS 1.27MM"

# Update the selected cell to mirror the saved selection in the workbook
$ws.Range("C18").Select()
